$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SamplePojo")

# New "sampleEnum" column header
$ws.Range("H1").Value = "sampleEnum"

# Enum values for each data row
$ws.Range("H2").Value = "SAMPLE1"
$ws.Range("H3").Value = "SAMPLE2"
$ws.Range("H4").Value = "SAMPLE1"
$ws.Range("H5").Value = "SAMPLE2"
$ws.Range("H6").Value = "SAMPLE1"

# Update selection to match the newly added column
$ws.Range("H1:H6").Select()
